# Applies the "Read data from excel" change: adds 14 vehicle/insurance fields
# (region..NCD) as new columns N:AA, with header row + one data row, matching
# the column-specific number formats (text for the short numeric-looking codes,
# thousands-separator for the insured amount) used by the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Claim the "#,##0" number-format style before any "@" (text) style cell is
# touched, so the style table is built in the same order as the source file
# (cellXfs[2] = #,##0, cellXfs[3] = @).
$ws.Cells.Item(2, 22).NumberFormat = "#,##0"

# --- Row 1 : new headers N1:AA1 ---
$ws.Cells.Item(1, 14).Value = "region"
$ws.Cells.Item(1, 15).Value = "car_number"
$ws.Cells.Item(1, 16).Value = "vehicle_type"
$ws.Cells.Item(1, 17).Value = "vehicle_make"
$ws.Cells.Item(1, 18).Value = "vehicle_model"
$ws.Cells.Item(1, 19).NumberFormat = "@"
$ws.Cells.Item(1, 19).Value = "seat_capacity"
$ws.Cells.Item(1, 20).NumberFormat = "@"
$ws.Cells.Item(1, 20).Value = "YOM"
$ws.Cells.Item(1, 21).Value = "vehicle_usage"
$ws.Cells.Item(1, 22).Value = "insured_amount"
$ws.Cells.Item(1, 23).NumberFormat = "@"
$ws.Cells.Item(1, 23).Value = "driving_exp"
$ws.Cells.Item(1, 24).Value = "garage"
$ws.Cells.Item(1, 25).Value = "package_type"
$ws.Cells.Item(1, 26).NumberFormat = "@"
$ws.Cells.Item(1, 26).Value = "voluntary_excess"
$ws.Cells.Item(1, 27).NumberFormat = "@"
$ws.Cells.Item(1, 27).Value = "NCD"

# --- Row 2 : new data values N2:AA2 (plus nic/F2 correction) ---
$ws.Cells.Item(2, 14).Value = "WP"
$ws.Cells.Item(2, 15).Value = "KR-9691"
$ws.Cells.Item(2, 16).Value = "Passenger Car"
$ws.Cells.Item(2, 17).Value = "PERODUA"
$ws.Cells.Item(2, 18).Value = "AMIZHR"
$ws.Cells.Item(2, 24).Value = "Standard"
$ws.Cells.Item(2, 25).Value = "Allianz Standard Package"
$ws.Cells.Item(2, 6).Value = "867361920V"
$ws.Cells.Item(2, 21).Value = "Hiring"
$ws.Cells.Item(2, 22).Value = 2300000
$ws.Cells.Item(2, 19).NumberFormat = "@"
$ws.Cells.Item(2, 19).Value = "4"
$ws.Cells.Item(2, 20).NumberFormat = "@"
$ws.Cells.Item(2, 20).Value = "2004"
$ws.Cells.Item(2, 23).NumberFormat = "@"
$ws.Cells.Item(2, 23).Value = "9"
$ws.Cells.Item(2, 26).NumberFormat = "@"
$ws.Cells.Item(2, 26).Value = "2000"
$ws.Cells.Item(2, 27).NumberFormat = "@"
$ws.Cells.Item(2, 27).Value = "60%"

# --- Column width (package_type is wider than the sheet default) ---
$ws.Columns("Y:Y").ColumnWidth = 29

# --- Final selection / active cell ---
$ws.Range("AA2").Select()
